# Insert a new weekly record as row 319 on the "Choclo" sheet, pushing the
# existing row 319 (and everything after it) down by one row. This mirrors
# the usual "add this week's price report" edit used for these weekly
# consolidated price sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 319; Excel shifts rows
# 319..419 down to 320..420 and the sheet dimension grows to A1:R420.
$ws.Rows("319:319").Insert()

# Populate the newly inserted row with the new price-report record.
$ws.Range("A319").Value = 10
$ws.Range("B319").Value = "Vega Modelo de Temuco"
$ws.Range("C319").Value = "La Araucanía"
$ws.Range("D319").Value = 44627
$ws.Range("E319").Value = 9
$ws.Range("F319").Value = 100112024
$ws.Range("G319").Value = "Choclo"
$ws.Range("H319").Value = "Choclero"
$ws.Range("I319").Value = "Primera"
$ws.Range("J319").Value = 8000
$ws.Range("K319").Value = 180
$ws.Range("L319").Value = 200
$ws.Range("M319").Value = 192
$ws.Range("N319").Value = "$/unidad"
$ws.Range("O319").Value = "Región del Maule"
$ws.Range("P319").Value = 192
$ws.Range("Q319").Value = 1
$ws.Range("R319").Value = "Hortaliza"
